$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 (被 / passive / 虛詞), shifting rows 9-13 up to become 8-12
$ws.Rows.Item(8).Delete()

# Match the author's final selection left in the saved file
[void]$ws.Range("F10").Select()
